# Auto-generated edit script applying the Tonberry_Profits update
# (scheduled market-price refresh across all job sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 16666989
$ws.Range("I92").Value = 20833612
$ws.Range("K92").Value = 20833612
$ws.Range("M92").Value = -20832364

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1148.4286
$ws.Range("J129").Value = 1200
$ws.Range("L129").Value = 3600
$ws.Range("N129").Value = -13600

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1202.2084
$ws.Range("I132").Value = 1113.3889
$ws.Range("K132").Value = 3340.1667
$ws.Range("M132").Value = -810.1666999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4006759.2
$ws.Range("I141").Value = 5603418.5
$ws.Range("K141").Value = 16810255.5
$ws.Range("M141").Value = -16805075.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 642.2
$ws.Range("J4").Value = 966.6667
$ws.Range("L4").Value = 966.6667
$ws.Range("N4").Value = -1198.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4884.778
$ws.Range("I32").Value = 4007.5813
$ws.Range("K32").Value = 4007.5813
$ws.Range("M32").Value = -3720.5813

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1709.5454
$ws.Range("J45").Value = 1854.7778
$ws.Range("L45").Value = 1854.7778
$ws.Range("N45").Value = -2608.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6051.5
$ws.Range("I61").Value = 6361.1177
$ws.Range("J61").Value = 5299.5713
$ws.Range("K61").Value = 6361.1177
$ws.Range("L61").Value = 5299.5713
$ws.Range("M61").Value = -6149.1177
$ws.Range("N61").Value = -5723.5713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1314.6086
$ws.Range("I74").Value = 457.72223
$ws.Range("K74").Value = 457.72223
$ws.Range("M74").Value = 416.27777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1314.6086
$ws.Range("I77").Value = 457.72223
$ws.Range("K77").Value = 2288.61115
$ws.Range("M77").Value = 2079.38885

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1958.8462
$ws.Range("I132").Value = 1544.4286
$ws.Range("J132").Value = 3699.4
$ws.Range("K132").Value = 4633.2858
$ws.Range("L132").Value = 11098.2
$ws.Range("M132").Value = -2103.2858
$ws.Range("N132").Value = -16158.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6051.5
$ws.Range("I136").Value = 6361.1177
$ws.Range("J136").Value = 5299.5713
$ws.Range("K136").Value = 19083.3531
$ws.Range("L136").Value = 15898.7139
$ws.Range("M136").Value = -16533.3531
$ws.Range("N136").Value = -20998.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 885.2857
$ws.Range("I22").Value = 780.6667
$ws.Range("J22").Value = 963.75
$ws.Range("K22").Value = 780.6667
$ws.Range("L22").Value = 963.75
$ws.Range("M22").Value = -607.6667
$ws.Range("N22").Value = -1309.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 22500
$ws.Range("I35").Value = 10000
$ws.Range("K35").Value = 10000
$ws.Range("M35").Value = -9690

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2675.7778
$ws.Range("I105").Value = 2322.75
$ws.Range("K105").Value = 2322.75
$ws.Range("M105").Value = -575.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4396.5835
$ws.Range("I134").Value = 4890.75
$ws.Range("J134").Value = 2667
$ws.Range("K134").Value = 14672.25
$ws.Range("L134").Value = 8001
$ws.Range("M134").Value = -12137.25
$ws.Range("N134").Value = -13071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5437584
$ws.Range("I58").Value = 14493557
$ws.Range("J58").Value = 3999.8
$ws.Range("K58").Value = 14493557
$ws.Range("L58").Value = 3999.8
$ws.Range("M58").Value = -14493354
$ws.Range("N58").Value = -4405.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3261.6155
$ws.Range("I132").Value = 1913.2858
$ws.Range("K132").Value = 5739.857400000001
$ws.Range("M132").Value = -3209.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5437584
$ws.Range("I136").Value = 14493557
$ws.Range("J136").Value = 3999.8
$ws.Range("K136").Value = 43480671
$ws.Range("L136").Value = 11999.4
$ws.Range("M136").Value = -43478121
$ws.Range("N136").Value = -17099.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 691.5
$ws.Range("I34").Value = 77.5
$ws.Range("J34").Value = 1305.5
$ws.Range("K34").Value = 232.5
$ws.Range("L34").Value = 3916.5
$ws.Range("M34").Value = -148.5
$ws.Range("N34").Value = -4084.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1954
$ws.Range("I64").Value = 412
$ws.Range("K64").Value = 1236
$ws.Range("M64").Value = -966

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 1954
$ws.Range("I67").Value = 412
$ws.Range("K67").Value = 1236
$ws.Range("M67").Value = -300

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 17406.771
$ws.Range("J131").Value = 19357.326
$ws.Range("L131").Value = 58071.978
$ws.Range("N131").Value = -68151.978

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3098.5789
$ws.Range("I138").Value = 1712.7273
$ws.Range("J138").Value = 5004.125
$ws.Range("K138").Value = 5138.1819
$ws.Range("L138").Value = 15012.375
$ws.Range("M138").Value = 1.818100000000413
$ws.Range("N138").Value = -25292.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3734.6667
$ws.Range("I70").Value = 3761
$ws.Range("J70").Value = 3701.75
$ws.Range("K70").Value = 3761
$ws.Range("L70").Value = 3701.75
$ws.Range("M70").Value = -3491
$ws.Range("N70").Value = -4241.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 3734.6667
$ws.Range("I73").Value = 3761
$ws.Range("J73").Value = 3701.75
$ws.Range("K73").Value = 3761
$ws.Range("L73").Value = 3701.75
$ws.Range("M73").Value = -2825
$ws.Range("N73").Value = -5573.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2735.875
$ws.Range("J122").Value = 2500
$ws.Range("L122").Value = 7500
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1480732.2
$ws.Range("I132").Value = 2025354.8
$ws.Range("K132").Value = 6076064.4
$ws.Range("M132").Value = -6073534.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1783.9767
$ws.Range("I132").Value = 1263.4445
$ws.Range("K132").Value = 3790.3335
$ws.Range("M132").Value = -1260.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2438.9
$ws.Range("I136").Value = 2339.5833
$ws.Range("J136").Value = 2587.875
$ws.Range("K136").Value = 7018.749899999999
$ws.Range("L136").Value = 7763.625
$ws.Range("M136").Value = -4468.749899999999
$ws.Range("N136").Value = -12863.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 29023.715
$ws.Range("I122").Value = 52879.734
$ws.Range("K122").Value = 158639.202
$ws.Range("M122").Value = -156189.202

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1981.6
$ws.Range("I126").Value = 1984.5454
$ws.Range("J126").Value = 1973.5
$ws.Range("K126").Value = 5953.6362
$ws.Range("L126").Value = 5920.5
$ws.Range("M126").Value = -3483.6362
$ws.Range("N126").Value = -10860.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1628.1277
$ws.Range("I132").Value = 1146.7188
$ws.Range("K132").Value = 3440.1564
$ws.Range("M132").Value = -910.1564000000003

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1095.1578
$ws.Range("I136").Value = 708.3929000000001
$ws.Range("K136").Value = 2125.1787
$ws.Range("M136").Value = 424.8212999999996

